# Fix multiple bugs in config file and program
# - Row 12 ("name%=%殡仪馆;category_code=142900 || category_code=208099") needed
#   parentheses around the OR'd category_code clause so it combines correctly
#   with the leading "name%=%..." condition:
#     name%=%殡仪馆;category_code=142900 || category_code=208099
#   becomes
#     name%=%殡仪馆;(category_code=142900 || category_code=208099)
# - That corrected row is highlighted with a yellow fill so it stands out.
# - The worksheet selection/scroll position is moved to show the corrected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C_CommunityArea.conf")

# --- Fix the formula/condition text in A12 -------------------------------
$ws.Range("A12").Value2 = "name%=%殡仪馆;(category_code=142900 || category_code=208099)"

# --- Highlight the corrected row with a solid yellow fill -----------------
$ws.Range("A12:C12").Interior.Color = 65535

# --- Update the view so the corrected row is selected/visible -------------
$ws.Range("A12:C12").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
